$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.468.46'
$ws.Range("E2").Value = '  -1.44%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.851.81'
$ws.Range("E3").Value = '  -0.47%  '

# Row 4
$ws.Range("E4").Value = '  +0.13%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '242.34'
$ws.Range("E5").Value = '  -0.88%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6306'
$ws.Range("E6").Value = '  -4.13%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.001'
$ws.Range("E7").Value = '  +0.13%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07564'
$ws.Range("E8").Value = '  -0.19%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2982'
$ws.Range("E9").Value = '  -0.48%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.48'
$ws.Range("E10").Value = '  -0.85%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07708'
$ws.Range("E11").Value = '  +0.83%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.922.39'
$ws.Range("E12").Value = '  +3.10%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.013'
$ws.Range("E13").Value = '  -1.35%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6901'
$ws.Range("E14").Value = '  -0.23%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '83.48'
$ws.Range("E15").Value = '  -0.59%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000009869'
$ws.Range("E16").Value = '  +1.53%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.170.12'
$ws.Range("E17").Value = '  +2.22%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.207'
$ws.Range("E18").Value = '  +0.94%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '29.624.85'
$ws.Range("E19").Value = '  -0.97%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '233.72'
$ws.Range("E20").Value = '  -1.62%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.56'
$ws.Range("E21").Value = '  -1.28%  '

# Row 22
$ws.Range("B22").Value = 'Chainlink'
$ws.Range("C22").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.716'
$ws.Range("E22").Value = '  -1.46%  '

# Row 23
$ws.Range("B23").Value = 'Dai'
$ws.Range("C23").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.002'
$ws.Range("E23").Value = '  +0.23%  '

# Row 24
$ws.Range("E24").Value = '  +0.13%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '155.51'
$ws.Range("E25").Value = '  -1.91%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1396'
$ws.Range("E26").Value = '  -3.16%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.485'
$ws.Range("E27").Value = '  -1.46%  '

# Row 28
$ws.Range("E28").Value = '  -1.06%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.478'
$ws.Range("E29").Value = '  -1.31%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05777'
$ws.Range("E30").Value = '  -4.39%  '

# Row 31
$ws.Range("E31").Value = '  -2.54%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.132'
$ws.Range("E32").Value = '  -0.94%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.028'
$ws.Range("E33").Value = '  -1.66%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.921'
$ws.Range("E34").Value = '  +2.43%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.170'
$ws.Range("E35").Value = '  -1.09%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7243'
$ws.Range("E36").Value = '  -1.40%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.593'
$ws.Range("E37").Value = '  -0.67%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.260.04'
$ws.Range("E38").Value = '  +4.28%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.806'
$ws.Range("E39").Value = '  -0.29%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01805'
$ws.Range("E40").Value = '  +0.32%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9058'
$ws.Range("E41").Value = '  -1.11%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.140'
$ws.Range("E42").Value = '  -2.77%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.096.86'
$ws.Range("E43").Value = '  +3.10%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.001'
$ws.Range("E44").Value = '  +0.07%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '68.02'
$ws.Range("E45").Value = '  +0.70%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '101.54'
$ws.Range("E46").Value = '  -0.13%  '

# Row 47
$ws.Range("B47").Value = 'BabyDogeCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000126'
$ws.Range("E47").Value = '  +4.12%  '

# Row 48
$ws.Range("B48").Value = 'Aptos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.362'
$ws.Range("E48").Value = '  -1.30%  '

# Row 49
$ws.Range("B49").Value = 'TheSandbox'
$ws.Range("C49").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.4037'
$ws.Range("E49").Value = '  -1.04%  '

# Row 50
$ws.Range("B50").Value = 'RenderToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.713'
$ws.Range("E50").Value = '  +1.18%  '

# Row 51
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '9.139'
$ws.Range("E51").Value = '  -0.90%  '
